$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "number included in dataset" (column E) counts per the refreshed
# literature-review tallies.
$ws.Range("E2").Value = 88
$ws.Range("E3").Value = 75
$ws.Range("E4").Value = 67
$ws.Range("E5").Value = 45
$ws.Range("E6").Value = 44
$ws.Range("E7").Value = 37
$ws.Range("E8").Value = 38
$ws.Range("E9").Value = 36
$ws.Range("E10").Value = 18
$ws.Range("E11").Value = 17
$ws.Range("E12").Value = 13
$ws.Range("E15").Value = 9
$ws.Range("E17").Value = 7
$ws.Range("E18").Value = 6
$ws.Range("E19").Value = 5
$ws.Range("E21").Value = 4
$ws.Range("E22").Value = 2
$ws.Range("E24").Value = 2
$ws.Range("E35").Value = 1
$ws.Range("E38").Value = 0

# Re-enter the "total" (D) and "dropped" (F) formulas across the full data
# range, in the same 64-row batches Excel used, so they recalc against the
# refreshed E values and regroup into shared formulas.
$ws.Range("D2:D65").Formula = "=SUM(B2, C2)"
$ws.Range("F2:F65").Formula = "=(D2-E2)"

$ws.Range("D66:D129").Formula = "=SUM(B66, C66)"
$ws.Range("F66:F129").Formula = "=(D66-E66)"

$ws.Range("D130:D193").Formula = "=SUM(B130, C130)"
$ws.Range("F130:F193").Formula = "=(D130-E130)"

$ws.Range("D194:D233").Formula = "=SUM(B194, C194)"
$ws.Range("F194:F233").Formula = "=(D194-E194)"
